# Applies the "addition of random transforms and multiprocess of grid mapping"
# commit to the KDE sheet of trainings.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KDE")

# ---------------------------------------------------------------------------
# 1) Insert a new row at 9 (old rows 9-12 shift down to 10-13). Excel copies
#    row 8's formatting into the freshly inserted row, which already gives us
#    the H9/K9 "quote-prefix" placeholder styles (s=6 / s=8) we need.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Insert()

# Row 10 (previously-blank old row 9) needs the same H/K placeholder styles
# that row 6 currently has (still untouched at this point in the script) -
# clone them across via a format-only paste.
$ws.Cells.Item(6,8).Copy()
$ws.Cells.Item(10,8).PasteSpecial(-4122)
$ws.Cells.Item(6,11).Copy()
$ws.Cells.Item(10,11).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Populate the new placeholder rows (9 and 10) with the experiment
#    parameters that used to live in rows 8 and 6 respectively.
# ---------------------------------------------------------------------------
$ws.Cells.Item(9,1).Value2 = 20
$ws.Cells.Item(9,2).Value2 = "yes"
$ws.Cells.Item(9,3).Value2 = "yes"
$ws.Cells.Item(9,4).Value2 = "no"
$ws.Cells.Item(9,5).Value2 = 64
$ws.Cells.Item(9,6).Value2 = 2
$ws.Cells.Item(9,7).Value2 = 0.4
$ws.Cells.Item(9,13).Value2 = 8
$ws.Cells.Item(9,14).Value2 = 5200

$ws.Cells.Item(10,1).Value2 = 20
$ws.Cells.Item(10,2).Value2 = "yes"
$ws.Cells.Item(10,3).Value2 = "yes"
$ws.Cells.Item(10,4).Value2 = "no"
$ws.Cells.Item(10,5).Value2 = 128
$ws.Cells.Item(10,6).Value2 = 2
$ws.Cells.Item(10,7).Value2 = 0.2
$ws.Cells.Item(10,13).Value2 = 8
$ws.Cells.Item(10,14).Value2 = 5200

# ---------------------------------------------------------------------------
# 3) Row 6: "adding whitening" run (50 epochs, kernel size 2, grid 64,
#    globavg model) - fill in the completed results.
# ---------------------------------------------------------------------------
$ws.Cells.Item(6,1).Value2 = 50
$ws.Cells.Item(6,5).Value2 = 64
$ws.Cells.Item(6,8).Value2 = 0.84406565656565602
$ws.Cells.Item(6,9).Value2 = 0.82632226185517699
$ws.Cells.Item(6,10).Value2 = 0.060786124619871602
$ws.Cells.Item(6,11).Value2 = 0.881126455456268
$ws.Cells.Item(6,12).Value2 = 0.056425763218908603
$ws.Cells.Item(6,13).Value2 = 12

# Re-apply the number formats (fixes the "quote prefix" style drift that
# happens when Excel sees a numeric value land in a style tagged
# quotePrefix="1").
$ws.Cells.Item(7,8).Copy()
$ws.Cells.Item(6,8).PasteSpecial(-4122)
$ws.Cells.Item(7,9).Copy()
$ws.Cells.Item(6,9).PasteSpecial(-4122)
$ws.Cells.Item(7,10).Copy()
$ws.Cells.Item(6,10).PasteSpecial(-4122)
$ws.Cells.Item(7,11).Copy()
$ws.Cells.Item(6,11).PasteSpecial(-4122)
$ws.Cells.Item(7,12).Copy()
$ws.Cells.Item(6,12).PasteSpecial(-4122)

$ws.Cells.Item(6,15).Value2 = "6:11:8"
$ws.Cells.Item(6,16).Value2 = "model_globavg"
$ws.Cells.Item(6,18).Value2 = "adding whitening"

$ws.Hyperlinks.Add($ws.Range("Q6"), "results\KDE\wl_da_tt_gd=64_ks=2_epoch=100_5200_globavg_whitening")
$ws.Cells.Item(2,17).Copy()
$ws.Cells.Item(6,17).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Row 8: "trying a deeper model" run (50 epochs, kernel size 2, grid 64,
#    globavg_deep model) - fill in the completed results.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8,1).Value2 = 50
$ws.Cells.Item(8,7).Value2 = 0.2
$ws.Cells.Item(8,8).Value2 = 0.83901515151515105
$ws.Cells.Item(8,9).Value2 = 0.82467885959482901
$ws.Cells.Item(8,10).Value2 = 0.062624890281997497
$ws.Cells.Item(8,11).Value2 = 0.85174654752233903
$ws.Cells.Item(8,12).Value2 = 0.059145503535370098
$ws.Cells.Item(8,13).Value2 = 12

$ws.Cells.Item(7,8).Copy()
$ws.Cells.Item(8,8).PasteSpecial(-4122)
$ws.Cells.Item(7,9).Copy()
$ws.Cells.Item(8,9).PasteSpecial(-4122)
$ws.Cells.Item(7,10).Copy()
$ws.Cells.Item(8,10).PasteSpecial(-4122)
$ws.Cells.Item(7,11).Copy()
$ws.Cells.Item(8,11).PasteSpecial(-4122)
$ws.Cells.Item(7,12).Copy()
$ws.Cells.Item(8,12).PasteSpecial(-4122)

$ws.Cells.Item(8,15).Value2 = "7:20:4"
$ws.Cells.Item(8,16).Value2 = "model_globavg_deep"
$ws.Cells.Item(8,18).Value2 = "trying a deeper model"

$ws.Hyperlinks.Add($ws.Range("Q8"), "results\KDE\wl_da_tt_gd=64_ks=2_epoch=100_5200_globavgdeep_whitening")
$ws.Cells.Item(2,17).Copy()
$ws.Cells.Item(8,17).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) Column width tweaks on P/Q (both are "best fit" columns that grew to
#    accommodate the longer strings just added).
# ---------------------------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 19.6
$ws.Columns.Item(17).ColumnWidth = 65.1

# ---------------------------------------------------------------------------
# 6) Selection cursor parked the way the source workbook left it.
# ---------------------------------------------------------------------------
$ws.Range("Q18").Select()
